$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the column headers in row 1: "_old" suffix -> "_FV2404",
#    "_new" suffix -> "_FV2410". The "diff" header (column K) is unchanged.
$fv2404Headers = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)
$fv2410Headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

for ($i = 0; $i -lt $fv2404Headers.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value2 = $fv2404Headers[$i]
}

for ($i = 0; $i -lt $fv2410Headers.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value2 = $fv2410Headers[$i]
}

# 2) Turn the used range into an actual Excel Table ("Table1") so the
#    header row gets the table/autofilter treatment.
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U66"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# 3) Freeze the header row (split/freeze after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Output ("Table created: " + $tbl.Name + " over " + $tbl.Range.Address())
